# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
#
# Updates the MSME Country Indicators - South Sudan Summary sheet:
#  - "Enterprises density (per 1000 people)" row (row 11): Micro/SMEs/MSMEs
#    values are corrected from 0.7 / 0.1 / 0.7 to the more precise
#    0.66 / 0.07 / 0.74.
#  - "Enterprises (% of total)" row (row 12): Micro/MSMEs values are
#    corrected from 89.8 / 99.7 to 89.83 / 99.73 (SMEs value 9.9 unchanged).
#
# All of these figures are stored as text in the workbook (matching the
# original authoring), so each value is entered with a leading apostrophe
# to force text interpretation, and the cell style is then reset back to
# "Normal" so that no spurious quote-prefix / number formatting is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("B11") "0.66"
Set-TextValue $ws.Range("C11") "0.07"
Set-TextValue $ws.Range("D11") "0.74"

Set-TextValue $ws.Range("B12") "89.83"
Set-TextValue $ws.Range("D12") "99.73"
